# Auto-generated PowerShell Excel COM-interop script
# Applies per-cell numeric updates to the Masamune_Profits.xlsx sheets
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 3133.3333
$ws.Range("I43").Value = 3142.8572
$ws.Range("J43").Value = 3000
$ws.Range("K43").Value = 3142.8572
$ws.Range("L43").Value = 3000
$ws.Range("M43").Value = -3073.8572
$ws.Range("N43").Value = -3138
$ws.Range("H62").Value = 7179.4443
$ws.Range("I62").Value = 8545
$ws.Range("J62").Value = 2400
$ws.Range("K62").Value = 8545
$ws.Range("L62").Value = 2400
$ws.Range("M62").Value = -7921
$ws.Range("N62").Value = -3648
$ws.Range("H65").Value = 7179.4443
$ws.Range("I65").Value = 8545
$ws.Range("J65").Value = 2400
$ws.Range("K65").Value = 42725
$ws.Range("L65").Value = 12000
$ws.Range("M65").Value = -39605
$ws.Range("N65").Value = -18240
$ws.Range("H132").Value = 23099.889
$ws.Range("I132").Value = 2792.95
$ws.Range("J132").Value = 185555.4
$ws.Range("K132").Value = 8378.849999999999
$ws.Range("L132").Value = 556666.2
$ws.Range("M132").Value = -5848.849999999999
$ws.Range("N132").Value = -561726.2
$ws.Range("H133").Value = 51921.9
$ws.Range("J133").Value = 51921.9
$ws.Range("L133").Value = 51921.9
$ws.Range("N133").Value = -62041.9
$ws.Range("H134").Value = 42180
$ws.Range("J134").Value = 42180
$ws.Range("L134").Value = 42180
$ws.Range("N134").Value = -52320
$ws.Range("H136").Value = 45522.223
$ws.Range("J136").Value = 45522.223
$ws.Range("L136").Value = 45522.223
$ws.Range("N136").Value = -55722.223
$ws.Range("H139").Value = 36607.5
$ws.Range("J139").Value = 36607.5
$ws.Range("L139").Value = 36607.5
$ws.Range("N139").Value = -46887.5
$ws.Range("H140").Value = 21731.111
$ws.Range("J140").Value = 21731.111
$ws.Range("L140").Value = 21731.111
$ws.Range("N140").Value = -32091.111

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2767.6758
$ws.Range("I61").Value = 2216.25
$ws.Range("K61").Value = 2216.25
$ws.Range("M61").Value = -2004.25
$ws.Range("H97").Value = 1753.8462
$ws.Range("I97").Value = 1314.2858
$ws.Range("J97").Value = 2266.6667
$ws.Range("K97").Value = 1314.2858
$ws.Range("L97").Value = 2266.6667
$ws.Range("M97").Value = -818.2858000000001
$ws.Range("N97").Value = -3258.6667
$ws.Range("H102").Value = 24235.607
$ws.Range("I102").Value = 3935.5557
$ws.Range("K102").Value = 3935.5557
$ws.Range("M102").Value = -2313.5557
$ws.Range("H110").Value = 1518.0625
$ws.Range("I110").Value = 1534.7742
$ws.Range("J110").Value = 1000
$ws.Range("K110").Value = 1534.7742
$ws.Range("L110").Value = 1000
$ws.Range("M110").Value = 510.2257999999999
$ws.Range("N110").Value = -5090
$ws.Range("H122").Value = 1695
$ws.Range("I122").Value = 1645.6451
$ws.Range("J122").Value = 1950
$ws.Range("K122").Value = 4936.9353
$ws.Range("L122").Value = 5850
$ws.Range("M122").Value = -2486.9353
$ws.Range("N122").Value = -10750
$ws.Range("H124").Value = 39429
$ws.Range("J124").Value = 39429
$ws.Range("L124").Value = 39429
$ws.Range("N124").Value = -49249
$ws.Range("H132").Value = 2526.3845
$ws.Range("I132").Value = 1721.5
$ws.Range("J132").Value = 4136.154
$ws.Range("K132").Value = 5164.5
$ws.Range("L132").Value = 12408.462
$ws.Range("M132").Value = -2634.5
$ws.Range("N132").Value = -17468.462
$ws.Range("H136").Value = 2767.6758
$ws.Range("I136").Value = 2216.25
$ws.Range("K136").Value = 6648.75
$ws.Range("M136").Value = -4098.75

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 2014.0667
$ws.Range("I94").Value = 1878.7273
$ws.Range("J94").Value = 2386.25
$ws.Range("K94").Value = 1878.7273
$ws.Range("L94").Value = 2386.25
$ws.Range("M94").Value = -1427.7273
$ws.Range("N94").Value = -3288.25
$ws.Range("H99").Value = 2483.1936
$ws.Range("I99").Value = 2508.6956
$ws.Range("J99").Value = 2409.875
$ws.Range("K99").Value = 2508.6956
$ws.Range("L99").Value = 2409.875
$ws.Range("M99").Value = -1010.6956
$ws.Range("N99").Value = -5405.875
$ws.Range("H107").Value = 1430.5883
$ws.Range("I107").Value = 1236.1818
$ws.Range("J107").Value = 1787
$ws.Range("K107").Value = 1236.1818
$ws.Range("L107").Value = 1787
$ws.Range("M107").Value = 683.8181999999999
$ws.Range("N107").Value = -5627
$ws.Range("H134").Value = 2269.3914
$ws.Range("I134").Value = 1978.8948
$ws.Range("J134").Value = 3649.25
$ws.Range("K134").Value = 5936.6844
$ws.Range("L134").Value = 10947.75
$ws.Range("M134").Value = -3401.6844
$ws.Range("N134").Value = -16017.75
$ws.Range("H140").Value = 31392.666
$ws.Range("J140").Value = 31392.666
$ws.Range("L140").Value = 31392.666
$ws.Range("N140").Value = -41752.666

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H52").Value = 44000
$ws.Range("J52").Value = 44000
$ws.Range("L52").Value = 44000
$ws.Range("N52").Value = -44588
$ws.Range("H122").Value = 63929.473
$ws.Range("J122").Value = 885.6
$ws.Range("L122").Value = 2656.8
$ws.Range("N122").Value = -7556.8
$ws.Range("H132").Value = 76018.31
$ws.Range("I132").Value = 1667.1818
$ws.Range("J132").Value = 178251.12
$ws.Range("K132").Value = 5001.5454
$ws.Range("L132").Value = 534753.36
$ws.Range("M132").Value = -2471.5454
$ws.Range("N132").Value = -539813.36
$ws.Range("H134").Value = 388693.75
$ws.Range("I134").Value = 415401.2
$ws.Range("J134").Value = 237351.67
$ws.Range("K134").Value = 1246203.6
$ws.Range("L134").Value = 712055.01
$ws.Range("M134").Value = -1243668.6
$ws.Range("N134").Value = -717125.01

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H120").Value = 1506000
$ws.Range("I120").Value = 3000000
$ws.Range("K120").Value = 9000000
$ws.Range("M120").Value = -8995162
$ws.Range("H132").Value = 4195.5884
$ws.Range("J132").Value = 6630.5557
$ws.Range("L132").Value = 59675.0013
$ws.Range("N132").Value = -64735.0013

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 35717424
$ws.Range("I132").Value = 52634060
$ws.Range("J132").Value = 4532.8887
$ws.Range("K132").Value = 157902180
$ws.Range("L132").Value = 13598.6661
$ws.Range("M132").Value = -157899650
$ws.Range("N132").Value = -18658.6661

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H92").Value = 45291.75
$ws.Range("J92").Value = 45291.75
$ws.Range("L92").Value = 45291.75
$ws.Range("N92").Value = -50283.75
$ws.Range("H122").Value = 2191.9092
$ws.Range("I122").Value = 2234.5557
$ws.Range("J122").Value = 2000
$ws.Range("K122").Value = 6703.6671
$ws.Range("L122").Value = 6000
$ws.Range("M122").Value = -4253.6671
$ws.Range("N122").Value = -10900

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 347.5
$ws.Range("I100").Value = 347.5
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 695
$ws.Range("L100").Value = 0
$ws.Range("M100").Value = -154
$ws.Range("N100").ClearContents()
$ws.Range("H122").Value = 2041717.8
$ws.Range("I122").Value = 3175427.5
$ws.Range("J122").Value = 1040
$ws.Range("K122").Value = 9526282.5
$ws.Range("L122").Value = 3120
$ws.Range("M122").Value = -9523832.5
$ws.Range("N122").Value = -8020
$ws.Range("H123").Value = 46500
$ws.Range("J123").Value = 46500
$ws.Range("L123").Value = 46500
$ws.Range("N123").Value = -56300

$wb.Save()
Write-Host "Applied all Masamune_Profits updates."